{"js": "// WIP extend to HS1 revision\n// Insert three new note paragraphs after \"V\u00e9rifier que l'indice total...\"\n// and add a leading space run inside the trailing bookmark paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"V\u00e9rifier que l'indice ...\" paragraph (last populated paragraph\n// before the trailing bookmark paragraph) robustly by its text, rather than\n// by a hard-coded index.\nconst anchorText =\n  \"V\u00e9rifier que l\\u2019indice total est coh\\u00e9rent avec une agr\\u00e9gation \" +\n  \"des indices au niveau ISIC.\";\n\nlet anchorParagraph = null;\nfor (const p of items) {\n  if (p.text.trim() === anchorText) {\n    anchorParagraph = p;\n    break;\n  }\n}\nif (!anchorParagraph) {\n  throw new Error(\"Anchor paragraph ('V\u00e9rifier que l'indice...') not found\");\n}\n\n// Helper: wrap a <w:p>\u2026</w:p> inner body in a minimal flat-OPC package so it\n// can be fed to Range.insertOoxml (needed to get the exact <w:proofErr/>\n// spell-check markers that plain insertText() would never produce).\nfunction flatOpcParagraph(innerXml) {\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" +\n    innerXml +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// The three new paragraphs' run content (second one keeps the curly\n// apostrophe used in the source doc, matching the diff verbatim).\nconst para1Runs =\n  '<w:r><w:t xml:space=\"preserve\">Ajouter un total </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>manuf</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r><w:t>.</w:t></w:r>\";\n\nconst para2Runs =\n  \"<w:r><w:t>Remplacer par la valeur seuil des 5% pour les observations \" +\n  \"filtr\\u00e9es. Puis s\\u2019il reste des trous boucher avec la moyenne \" +\n  \"niveau SH 4 chiffres.</w:t></w:r>\";\n\nconst para3Runs =\n  \"<w:r><w:t>On opte pour le filtrage non pond\\u00e9r\\u00e9.</w:t></w:r>\";\n\n// Create three empty paragraphs right after the anchor (insertParagraph\n// mirrors Range.InsertParagraphAfter, so it never disturbs the anchor's own\n// content or the paragraph that follows), then fill each via insertOoxml so\n// the spell-check <w:proofErr/> markers come through untouched.\nconst newPara1 = anchorParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\nconst newPara2 = newPara1.insertParagraph(\"\", \"After\");\nawait context.sync();\nconst newPara3 = newPara2.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nnewPara1.getRange().insertOoxml(flatOpcParagraph(para1Runs), \"Replace\");\nawait context.sync();\nnewPara2.getRange().insertOoxml(flatOpcParagraph(para2Runs), \"Replace\");\nawait context.sync();\nnewPara3.getRange().insertOoxml(flatOpcParagraph(para3Runs), \"Replace\");\nawait context.sync();\n\n// Re-fetch paragraphs (indices/anchors shifted after the inserts above) and\n// prepend a lone space run to the final (bookmark) paragraph, right before\n// its existing <w:bookmarkStart/>.\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph =\n  refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\nlastParagraph.getRange(\"Start\").insertText(\" \", \"Before\");\nawait context.sync();\n", "ps1": "# WIP extend to HS1 revision\n# Insert three new note paragraphs after \"V\u00e9rifier que l'indice total...\"\n# and add a leading space run inside the trailing bookmark paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the \"V\u00e9rifier que l'indice total...\" paragraph robustly (ASCII-safe\n# wildcard match so we don't depend on exact curly-quote/accent byte forms).\n$anchor = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*indice total est*niveau ISIC*\") {\n        $anchor = $p\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph ('V\u00e9rifier que l'indice...') not found\"\n}\n\n# Create three empty paragraphs right after the anchor. InsertParagraphAfter\n# never disturbs the anchor's own content nor the paragraph that follows\n# (the trailing bookmark paragraph), so it keeps them intact.\n$anchor.Range.InsertParagraphAfter()\n$p1 = $anchor.Next()\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $p1.Next()\n\n$p2.Range.InsertParagraphAfter()\n$p3 = $p2.Next()\n\n# Fill each new paragraph via Range.InsertXML (flat-OPC wrapped) so the\n# <w:proofErr/> spell-check markers around \"manuf\" come through verbatim --\n# plain Range.Text assignment would never produce them.\n$xml1 = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r><w:t xml:space=\"preserve\">Ajouter un total </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>manuf</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$xml2 = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r><w:t>Remplacer par la valeur seuil des 5% pour les observations filtr\u00e9es. Puis s\u2019il reste des trous boucher avec la moyenne niveau SH 4 chiffres.</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$xml3 = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:r><w:t>On opte pour le filtrage non pond\u00e9r\u00e9.</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$p1.Range.InsertXML($xml1)\n$p2.Range.InsertXML($xml2)\n$p3.Range.InsertXML($xml3)\n\n# Prepend a lone space run to the final (bookmark) paragraph, right before\n# its existing bookmarkStart/bookmarkEnd.\n$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)\n$rLast = $pLast.Range\n$rLast.Collapse(1)  # wdCollapseStart\n$rLast.InsertBefore(\" \")\n"}
